$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.284.60'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '2.315.24'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.07'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.73'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.576'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('D9').Value = '2.313.37'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.55'
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.43'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').Value = '60.336.78'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').Value = '2.728.74'
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = '2.312.21'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.56'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '313.36'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  -3.44%  '
$ws.Range('E23').Value = '  +0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.08'
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.85'
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('E28').Value = '  +3.37%  '
$ws.Range('E29').Value = '  +8.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.13'
$ws.Range('E30').Value = '  +1.15%  '
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.95'
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.35'
$ws.Range('E35').Value = '  -3.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.96'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '316.96'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.01'
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '136.71'
$ws.Range('E43').Value = '  -3.79%  '
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0939'
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.05'
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.563'
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').Value = '0.0₆0216'
$ws.Range('E50').Value = '  +6.34%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.87'
$ws.Range('E51').Value = '  +0.18%  '
